$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Fill in the activity description and coffee count for the two newly
# journaled days (rows 23 and 24).
$ws.Range("C23").Value = "Programmation du firmware, lecture/ecriture carte SD, gestion du fichier de config."
$ws.Range("E23").Value = 3

$ws.Range("C24").Value = "Programmation du firmware, gestion carte SD, centrale inertielle et GNSS."
$ws.Range("E24").Value = 1

# Update the view: scroll so the top-left cell resets and select E25.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E25").Select()
